$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.992.20'
$ws.Range("E2").Value = '  +1.63%  '
$ws.Range("D3").Value = '1.861.26'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").Value = '335.63'
$ws.Range("E5").Value = '  +0.41%  '
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("E7").Value = '  +1.18%  '
$ws.Range("D8").Value = '0.3901'
$ws.Range("D9").Value = '46.75'
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("D10").Value = '0.07958'
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("D11").Value = '0.9815'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("D12").Value = '21.50'
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '5.936'
$ws.Range("E13").Value = '  -0.36%  '
$ws.Range("D14").Value = '1.838.02'
$ws.Range("E14").Value = '  -0.80%  '
$ws.Range("D15").Value = '7.191'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '91.41'
$ws.Range("E16").Value = '  +3.26%  '
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.52%  '
$ws.Range("E18").Value = '  +0.45%  '
$ws.Range("D19").Value = '0.06615'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").Value = '17.48'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  -0.31%  '
$ws.Range("D22").Value = '28.028.66'
$ws.Range("E22").Value = '  +1.76%  '
$ws.Range("D23").Value = '5.393'
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("D24").Value = '10.94'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = '2.286'
$ws.Range("E25").Value = '  -1.22%  '
$ws.Range("D26").Value = '159.26'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '2.065.47'
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").Value = '19.52'
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("D29").Value = '2.101'
$ws.Range("E29").Value = '  -0.65%  '
$ws.Range("D30").Value = '5.459'
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").Value = '119.26'
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("D32").Value = '0.9602'
$ws.Range("E32").Value = '  -1.64%  '
$ws.Range("D33").Value = '0.09472'
$ws.Range("E33").Value = '  +0.77%  '
$ws.Range("D34").Value = '3.577'
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").Value = '5.301'
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = '1.347'
$ws.Range("E36").Value = '  +0.60%  '
$ws.Range("D37").Value = '0.06079'
$ws.Range("E37").Value = '  +0.52%  '
$ws.Range("D38").Value = '0.02252'
$ws.Range("E38").Value = '  +1.17%  '
$ws.Range("D39").Value = '8.293'
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = '1.162'
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("D42").Value = '0.5926'
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("D43").Value = '0.1866'
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("D44").Value = '10.20'
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").Value = '1.279'
$ws.Range("E45").Value = '  +3.42%  '
$ws.Range("D46").Value = '0.5544'
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("D47").Value = '12.15'
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").Value = '1.950'
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("D49").Value = '0.06876'
$ws.Range("E49").Value = '  +2.61%  '
$ws.Range("D50").Value = '111.44'
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("D51").Value = '1.001'
